$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H1").Value2 = '<%=comment.create_usr_id_lbl%><%selectList.create_usr_id = data.findAllUsr.map((item) => item.lbl)%><%_dataValidation_({ sqref: `${ _col }2:${ _col }${ _lastRow }`, formula1: `"${ selectList.create_usr_id.join(",") }"` })%>'
$ws.Range("I1").Value2 = '<%=comment.create_time_lbl%>'
$ws.Range("J1").Value2 = '<%=comment.update_usr_id_lbl%><%selectList.update_usr_id = data.findAllUsr.map((item) => item.lbl)%><%_dataValidation_({ sqref: `${ _col }2:${ _col }${ _lastRow }`, formula1: `"${ selectList.update_usr_id.join(",") }"` })%>'
$ws.Range("K1").Value2 = '<%=comment.update_time_lbl%>'

$ws.Range("H2").Value2 = '<%=model.create_usr_id_lbl%>'
$ws.Range("I2").Value2 = '<%~model.create_time ? new Date(model.create_time) : ""%>'
$ws.Range("J2").Value2 = '<%=model.update_usr_id_lbl%>'
$ws.Range("K2").Value2 = '<%~model.update_time ? new Date(model.update_time) : ""%>'
